# Pouya Finance / Shapna (Esfahan Oil Refining) - Income Statement (USD)
# Rolls the 5-year trailing window forward by one fiscal year:
#   drops FY1396/12, adds FY1401/12, and refreshes every figure that shifts
#   with it (publish dates + all reported line items).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8: "12 ماهه منتهی به <year>/12" financial-period headers ----
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---- Row 9: "تاریخ انتشار" (publish date) headers ----
$ws.Range("D9").Value = "1399-04-19 (8)"
$ws.Range("E9").Value = "1400-04-26 (10)"
$ws.Range("F9").Value = "1401-04-21 (10)"
$ws.Range("G9").Value = "1402-02-29 (8)"
$ws.Range("H9").Value = "1402-02-29"

# ---- Data rows: each row's D:H values shift one column left, with a new
#      right-most (H) figure for FY1401/12. "-" cells stay textual "-". ----

# Row 11: فروش (Sales)
$ws.Range("D11").Value = 4798376
$ws.Range("E11").Value = 6753337
$ws.Range("F11").Value = 4907281
$ws.Range("G11").Value = 8075622
$ws.Range("H11").Value = 9627090

# Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold)
$ws.Range("D12").Value = -4249647
$ws.Range("E12").Value = -6326681
$ws.Range("F12").Value = -4340106
$ws.Range("G12").Value = -7340466
$ws.Range("H12").Value = -8367351

# Row 13: سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 548729
$ws.Range("E13").Value = 426657
$ws.Range("F13").Value = 567175
$ws.Range("G13").Value = 735156
$ws.Range("H13").Value = 1259740

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
$ws.Range("D14").Value = -51527
$ws.Range("E14").Value = -56616
$ws.Range("F14").Value = -36761
$ws.Range("G14").Value = -47825
$ws.Range("H14").Value = -74158

# Row 15: هزینه کاهش ارزش دریافتنی‌‏ها (stays "-" across the board)
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense, net)
$ws.Range("D16").Value = 92934
$ws.Range("E16").Value = -154982
$ws.Range("F16").Value = 171040
$ws.Range("G16").Value = 45421
$ws.Range("H16").Value = -628

# Row 17: سود (زیان) عملیاتی (Operating profit)
$ws.Range("D17").Value = 590136
$ws.Range("E17").Value = 215059
$ws.Range("F17").Value = 701455
$ws.Range("G17").Value = 732752
$ws.Range("H17").Value = 1184953

# Row 18: هزینه های مالی (Financial expenses) - D/E remain "-", F now numeric
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = -393
$ws.Range("G18").Value = -24806
$ws.Range("H18").Value = -46506

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense, net)
$ws.Range("D19").Value = 1631
$ws.Range("E19").Value = 32245
$ws.Range("F19").Value = 32199
$ws.Range("G19").Value = 30552
$ws.Range("H19").Value = 73347

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit from continuing ops)
$ws.Range("D20").Value = 591766
$ws.Range("E20").Value = 247304
$ws.Range("F20").Value = 733261
$ws.Range("G20").Value = 738497
$ws.Range("H20").Value = 1211794

# Row 21: مالیات (Tax)
$ws.Range("D21").Value = -89829
$ws.Range("E21").Value = -24945
$ws.Range("F21").Value = -76704
$ws.Range("G21").Value = -110951
$ws.Range("H21").Value = -129697

# Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing ops)
$ws.Range("D22").Value = 501937
$ws.Range("E22").Value = 222359
$ws.Range("F22").Value = 656558
$ws.Range("G22").Value = 627546
$ws.Range("H22").Value = 1082098

# Row 24: سود (زیان) خالص (Net profit) - mirrors row 22 (row 23 discontinued ops stays "-")
$ws.Range("D24").Value = 501937
$ws.Range("E24").Value = 222359
$ws.Range("F24").Value = 656558
$ws.Range("G24").Value = 627546
$ws.Range("H24").Value = 1082098

# Row 25: سود هر سهم پس از کسر مالیات (EPS after tax) - E25 becomes "-", F25 becomes 0
$ws.Range("E25").Value = "-"
$ws.Range("F25").Value = 0

# Row 26: سرمایه (Capital)
$ws.Range("D26").Value = 504166
$ws.Range("E26").Value = 397562
$ws.Range("F26").Value = 613214
$ws.Range("G26").Value = 720097
$ws.Range("H26").Value = 643259
